$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row1 = @(-0.0091819980337988858, -0.0053427473055305333, -0.0097178397593929607, -0.0078811667531655839, -0.0067634929278475508, -0.0049932384171355758, -0.008202100189464212, -0.0066603716753310619, -0.0089338388189316546, -0.0073568293760257119, -0.006364244560287495, -0.0070290918241296303, -0.0039412687495760648, -0.007535339668414495, -0.0070801204402612629, -0.009078543113916468, -0.0064771885594827432, -0.0081219678240546829, -0.0066796883322349414, -0.0089110120157371801, -0.0049066621457328529, -0.011526503266864021, -0.0073948240892507391, -0.0071967720101541871, -0.0082402870094704806, -0.0053815088750023997, -0.0098842157871807269, -0.0073633932725985958, -0.0040452998982180629, -0.003779290335538994, -0.0081270218874019036, -0.0062317465081067301, -0.0068178820700799392, -0.0058416435299399267, -0.0068373014143120669, -0.0058529860906330689, -0.0080135237588093589, -0.0051537040383196923, -0.0087731631659111933, -0.0071601416373648181, -0.0053873455088755259, -0.012617296899878113, -0.010527547746610699, -0.0090165385053548857, -0.0076871469556173411, -0.0024652727631233927, -0.0088851882371918552, -0.0075144473046717692, -0.0081884243438038323, -0.0063042396540661455, -0.0088041415151515878, -0.0095775892894473733, -0.0055972871694003196, -0.002306988126884936, -0.0097270025403645745, -0.0090050599371545755, -0.0076190049684183537, -0.0064638461647808102, -0.015225631402538749, -0.0084554921344091906, -0.0098629206864996998, -0.0047493957300026356, -0.0092259643470953001, -0.0064654094262127223, -0.0092147484819883208, -0.0066679697818575598, -0.0077214859741516111, -0.010972375493972035, -0.0064086852724880112, -0.011760800226973404, -0.0043846211069583224, -0.010656614296774614, -0.0086884253004113621, -0.0073866905380045389, -0.0077271926161336614, -0.005860950243410541, -0.0069167299166315623, -0.0074879807233557113, -0.0078462833340510879, -0.0073275763890100157, -0.0074196567523365604, -0.011165517553005906, -0.009168227243620446, -0.0071970781669501568, -0.0062599346931555642, -0.0049141422401480532, -0.0024983407259549658, -0.0038164967729898526, -0.0060727996022276031, -0.0091602360170399497, -0.007800838331421223, -0.0068135758894969194, -0.0069235678183279317, -0.0055522896140455481, -0.011152584880001815, -0.0077244956960338196, -0.011409432576485841, -0.011607864859644495, -0.009361438685423825, -0.005814113013069026)
$row2 = @(-0.013668417412427206, -0.0078267374716220906, -0.014636608781829567, -0.011860560747251581, -0.010078562048807719, -0.0074730821472695219, -0.012359374437949757, -0.010195109174677589, -0.01311038760724997, -0.011165730696261289, -0.0094796287817438055, -0.010387886484025928, -0.0059463341965112718, -0.011492876578089006, -0.010665459417696983, -0.01334587155100418, -0.0097638671238293694, -0.012233293402251389, -0.010050524680997052, -0.013680628300326328, -0.0073926353017200004, -0.017288905067571525, -0.011197785013626897, -0.010833677789937877, -0.012195259054993152, -0.0079193602482837272, -0.014662512120028635, -0.010958832806498403, -0.0059364004057332618, -0.0056995681222105988, -0.012091501289667269, -0.0092553564095355628, -0.010225389540608847, -0.0086628359503298728, -0.010205717443401449, -0.0085908168672839264, -0.012236674802189299, -0.0078357542125067452, -0.013140854780175973, -0.01069955688594585, -0.008173219029415103, -0.019120134123695409, -0.015872932106917821, -0.013496771671036849, -0.011460572792562225, -0.0036187591997602041, -0.013198640637618926, -0.011314161647484075, -0.012187979755763263, -0.009304559092095592, -0.013350539434589791, -0.014351474983120296, -0.0083670199823945423, -0.003449868327568344, -0.014534714038238571, -0.013389361486656847, -0.01139281281248454, -0.0097734237072751896, -0.022913431367783836, -0.012479696075761817, -0.014771477561376508, -0.0069758736849603971, -0.013849988779386391, -0.0098780041911754497, -0.013597119386565893, -0.0098984443535970416, -0.011424662845517778, -0.016409369831070497, -0.0095916667088482514, -0.017690965321226756, -0.0063931566264119289, -0.015741547925140306, -0.013197083956655939, -0.011171606567734721, -0.011791023671404327, -0.0086625452184432381, -0.010300838661402828, -0.011008632488365528, -0.011592560984078331, -0.01097163731454757, -0.011231209284413846, -0.016812363053510221, -0.013802611464357373, -0.010603097557410556, -0.009266412788590599, -0.0075346301276741391, -0.0036646699672830741, -0.0057489364282558219, -0.0092662518291571436, -0.013851755063253972, -0.011716678785369595, -0.010090717158488667, -0.010397354463191831, -0.0085716205973731148, -0.016760787766782432, -0.011402986761400032, -0.017116228249614014, -0.017400637640131982, -0.013883345012292382, -0.0087124838674137693)

for ($i = 0; $i -lt $row1.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $row1[$i]
    $ws.Cells.Item(2, $i + 1).Value = $row2[$i]
}
